# New entries to data base
# - Rename a few "focus" terms (G column) to use "Near-IR" / "Mid-IR" wording
# - Add three new package rows (spectacles, baseline, geoSpectral)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing focus terms -------------------------------------------------
$ws.Range("G15").Value = "Vis-Near-IR"
$ws.Range("G25").Value = "Time series (Near-IR)"
$ws.Range("G46").Value = "Near-IR, Mid-IR"

# --- New row 55: spectacles --------------------------------------------------------
$ws.Range("A55").Value = "spectacles"
$ws.Range("B55").Value = "https://cran.r-project.org/package=spectacles"
$ws.Range("C55").Value = "https://github.com/pierreroudier/spectacles/"
$ws.Range("E55").Value = "R"
$ws.Range("F55").Value = "Storing, manipulation and analysis "

# --- New row 56: baseline -----------------------------------------------------------
$ws.Range("A56").Value = "baseline"
$ws.Range("B56").Value = "https://cran.r-project.org/package=baseline"
$ws.Range("E56").Value = "R"
$ws.Range("F56").Value = "Baseline correction"

# --- New row 57: geoSpectral ---------------------------------------------------------
$ws.Range("A57").Value = "geoSpectral"
$ws.Range("B57").Value = "https://cran.r-project.org/package=baseline"
$ws.Range("E57").Value = "R"
$ws.Range("F57").Value = "Workflow for data sets with space/time/spectral dimensions"
$ws.Range("G57").Value = "Mid-Ir"

# --- Update view / selection state to match author's final window state --------------
$ws.Range("F6").Select() | Out-Null
